# Apply the "fix issue, and add some issue" commit to the exchange example sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data fixes --------------------------------------------------------
# Row 5's "To" address was wrong - point it at the correct recipient address.
$ws.Range("A5").Value = "0xf881a94423f22ee9a0e3e1442f515f43c966b7ed"

# Amount column was off by three orders of magnitude.
$ws.Range("B3").Value = 1.88
$ws.Range("B4").Value = 1.88
$ws.Range("B5").Value = 1.88

# ChainTag updated.
$ws.Range("C3").Value = "0x27"
$ws.Range("C4").Value = "0x27"
$ws.Range("C5").Value = "0x27"

# BlockRef updated.
$ws.Range("D3").Value = "0x0000695540f491a5"
$ws.Range("D4").Value = "0x0000695540f491a5"
$ws.Range("D5").Value = "0x0000695540f491a5"

# --- Formatting fix ------------------------------------------------------
# The BlockRef column (D) used to carry its own one-off style; align all
# three rows with the first ChainTag cell's formatting it was meant to match.
$ws.Range("C3").Copy()
$ws.Range("D3:D5").PasteSpecial(-4122)
$wb.Application.CutCopyMode = $false

# --- Leave the selection where the author last left it -------------------
$ws.Range("A8").Select()
